$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 4970.8
$ws.Range("I29").Value = 3914.3333
$ws.Range("J29").Value = 6555.5
$ws.Range("K29").Value = 11742.9999
$ws.Range("L29").Value = 19666.5
$ws.Range("M29").Value = -11461.9999
$ws.Range("N29").Value = -20228.5
# Row 38
$ws.Range("H38").Value = 6916.5
$ws.Range("I38").Value = 5833
$ws.Range("J38").Value = 8000
$ws.Range("K38").Value = 17499
$ws.Range("L38").Value = 24000
$ws.Range("M38").Value = -17127
$ws.Range("N38").Value = -24744
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
# Row 51
$ws.Range("H51").Value = 2085
$ws.Range("I51").Value = 2085
$ws.Range("K51").Value = 2085
$ws.Range("M51").Value = -1601
# Row 58
$ws.Range("H58").Value = 2063.6667
$ws.Range("J58").Value = 2999.25
$ws.Range("L58").Value = 8997.75
$ws.Range("N58").Value = -9297.75
# Row 86
$ws.Range("H86").Value = 1673.2727
$ws.Range("I86").Value = 1673.2727
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1673.2727
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -550.2727
# Row 89
$ws.Range("H89").Value = 1673.2727
$ws.Range("I89").Value = 1673.2727
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 8366.363499999999
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2750.363499999999
$ws.Range("N89").ClearContents()
# Row 125
$ws.Range("H125").Value = 5606136
$ws.Range("I125").Value = 440.25
$ws.Range("J125").Value = 9343266
$ws.Range("K125").Value = 3962.25
$ws.Range("L125").Value = 84089394
$ws.Range("M125").Value = -1502.25
$ws.Range("N125").Value = -84094314

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 1185.4
$ws.Range("I122").Value = 808
$ws.Range("J122").Value = 1751.5
$ws.Range("K122").Value = 2424
$ws.Range("L122").Value = 5254.5
$ws.Range("M122").Value = 26
$ws.Range("N122").Value = -10154.5
# Row 133
$ws.Range("H133").Value = 35857
$ws.Range("J133").Value = 35857
$ws.Range("L133").Value = 35857
$ws.Range("N133").Value = -40917

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 676.3333
$ws.Range("I107").Value = 593.3889
$ws.Range("J107").Value = 842.2222
$ws.Range("K107").Value = 593.3889
$ws.Range("L107").Value = 842.2222
$ws.Range("M107").Value = 1326.6111
$ws.Range("N107").Value = -4682.2222

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 592.6
$ws.Range("I16").Value = 592.6
$ws.Range("K16").Value = 592.6
$ws.Range("M16").Value = -305.6
# Row 31
$ws.Range("H31").Value = 1874.4117
$ws.Range("I31").Value = 1445.909
$ws.Range("J31").Value = 2660
$ws.Range("K31").Value = 1445.909
$ws.Range("L31").Value = 2660
$ws.Range("M31").Value = -1150.909
$ws.Range("N31").Value = -3250
# Row 34
$ws.Range("H34").Value = 1874.4117
$ws.Range("I34").Value = 1445.909
$ws.Range("J34").Value = 2660
$ws.Range("K34").Value = 1445.909
$ws.Range("L34").Value = 2660
$ws.Range("M34").Value = -1243.909
$ws.Range("N34").Value = -3064
# Row 58
$ws.Range("H58").Value = 2313.6287
$ws.Range("I58").Value = 1216.9
$ws.Range("K58").Value = 1216.9
$ws.Range("M58").Value = -1013.9
# Row 99
$ws.Range("H99").Value = 8930089
$ws.Range("I99").Value = 8930089
$ws.Range("K99").Value = 8930089
$ws.Range("M99").Value = -8928591
# Row 113
$ws.Range("H113").Value = 592.6
$ws.Range("I113").Value = 592.6
$ws.Range("K113").Value = 592.6
$ws.Range("M113").Value = 1577.4
# Row 126
$ws.Range("H126").Value = 8930089
$ws.Range("I126").Value = 8930089
$ws.Range("K126").Value = 26790267
$ws.Range("M126").Value = -26787797
# Row 136
$ws.Range("H136").Value = 2313.6287
$ws.Range("I136").Value = 1216.9
$ws.Range("K136").Value = 3650.7
$ws.Range("M136").Value = -1100.7

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 21740102
$ws.Range("I113").Value = 569.8333
$ws.Range("J113").Value = 29412878
$ws.Range("K113").Value = 1709.4999
$ws.Range("L113").Value = 88238634
$ws.Range("M113").Value = 460.5001
$ws.Range("N113").Value = -88242974
# Row 122
$ws.Range("H122").Value = 1001
$ws.Range("I122").Value = 335.66666
$ws.Range("K122").Value = 3020.99994
$ws.Range("M122").Value = -570.9999399999997

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2800
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1702
$ws.Range("N80").Value = -4996
# Row 83
$ws.Range("H83").Value = 2800
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 13500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -8508
$ws.Range("N83").Value = -24984
# Row 113
$ws.Range("H113").Value = 1194.1666
$ws.Range("I113").Value = 941.25
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 941.25
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 1228.75
$ws.Range("N113").Value = -6040
# Row 126
$ws.Range("H126").Value = 2273.1765
$ws.Range("I126").Value = 1523.2667
$ws.Range("K126").Value = 4569.800099999999
$ws.Range("M126").Value = -2099.800099999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3323.5293
$ws.Range("I7").Value = 2333.3333
$ws.Range("J7").Value = 3535.7144
$ws.Range("K7").Value = 2333.3333
$ws.Range("L7").Value = 3535.7144
$ws.Range("M7").Value = -2221.3333
$ws.Range("N7").Value = -3759.7144
# Row 126
$ws.Range("H126").Value = 3323.5293
$ws.Range("I126").Value = 2333.3333
$ws.Range("J126").Value = 3535.7144
$ws.Range("K126").Value = 6999.999899999999
$ws.Range("L126").Value = 10607.1432
$ws.Range("M126").Value = -4529.999899999999
$ws.Range("N126").Value = -15547.1432

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 44469332
$ws.Range("J2").Value = 43499
$ws.Range("L2").Value = 43499
$ws.Range("N2").Value = -43723
# Row 126
$ws.Range("H126").Value = 33659.516
$ws.Range("I126").Value = 42962.293
$ws.Range("J126").Value = 1764.2858
$ws.Range("K126").Value = 128886.879
$ws.Range("L126").Value = 5292.857400000001
$ws.Range("M126").Value = -126416.879
$ws.Range("N126").Value = -10232.8574
